# Apply the "add correct results and results to debug of domapriori" edit:
#  - reorder the object lists inside several rule descriptions on the
#    "Reguły" sheet
#  - update the "coverage" column values for rules 1 and 3 on the
#    "Statystyki reguł" sheet

$wb = $excel.ActiveWorkbook

# --- Sheet "Reguły" (rule descriptions, column B) ---
$rules = $wb.Worksheets.Item(8)

$rules.Range("B2").Value = "(attempts >=  3.0) & (pregnancy <=  0.0) => (class <= 1) ['a3', 'a7', 'a1']"
$rules.Range("B3").Value = "(sperm >=  3.0) => (class <= 1) ['a22', 'a25']"
$rules.Range("B4").Value = "(age >=  40.0) & (pregnancy <=  0.0) => (class <= 1) ['a3', 'a15']"
$rules.Range("B6").Value = "(age >=  42.0) => (class <= 1) ['a14', 'a3']"
$rules.Range("B7").Value = "(age <=  31.0) & (attempts <=  1.0) & (endometrium <=  1.0) => (class >= 2) ['a24', 'a9', 'a11', 'a12']"
$rules.Range("B8").Value = "(frozen_embryos >=  8.0) & (sperm <=  1.0) => (class >= 2) ['a16', 'a6']"

# --- Sheet "Statystyki reguł" (coverage column, column C) ---
$stats = $wb.Worksheets.Item(9)

$stats.Range("C2").Value = 0.375
$stats.Range("C4").Value = 0.25
